$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: mark "carro para cama quente" as concluído, add note about
# missing stopper button.
$ws.Range("F31").Value = "concluído"
$ws.Range("G31").Value = "3 peças, falta botão stopper"

# Row 34 (new): parafuso CHC M5x47 for the cama quente carriage.
$ws.Range("A34").Value = "parafuso CHC M5x47"
$ws.Range("B34").Value = 4
$ws.Range("C34").Value = "comprar"
$ws.Range("E34").Value = "não"
$ws.Range("G34").Value = "para carro da cama quente"

# Row 35 (new): parafuso CHC M3x for the cama quente carriage.
$ws.Range("A35").Value = "parafuso CHC M3x"
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = "comprar"
$ws.Range("E35").Value = "não"
$ws.Range("G35").Value = "para carro da cama quente"

# Row 37: new formatted (underlined + bordered) empty cell at F37.
$ws.Range("F37").Font.Underline = $true
$ws.Range("F37").Borders.Item(7).LineStyle = 1

# Update the active view to match the edited area.
$ws.Range("F37").Select()
